$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Change 1: paragraph "Check if every information ... before implement
# controller" (numId=4) gains an eastAsia font hint on its paragraph mark
# (i.e. a <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> inside its <w:pPr>).
# ---------------------------------------------------------------------------
$checkPara = $d.Paragraphs.Item(2)
if ($checkPara.Range.Text.StartsWith("Check ")) {
    $checkXml = "<w:p $W w:rsidR='005A4A50' w:rsidRDefault='005A4A50' w:rsidP='005A4A50'>" +
        "<w:pPr>" +
          "<w:pStyle w:val='a7'/>" +
          "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr>" +
          "<w:ind w:firstLineChars='0'/>" +
          "<w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr>" +
        "</w:pPr>" +
        "<w:r w:rsidRPr='005A4A50'><w:t xml:space='preserve'>Check </w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>if </w:t></w:r>" +
        "<w:r w:rsidRPr='005A4A50'><w:t>every information from sensor and corresponding calculation</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'> is correct</w:t></w:r>" +
        "<w:r w:rsidRPr='005A4A50'><w:t xml:space='preserve'> before implement controller</w:t></w:r>" +
      "</w:p>"
    $checkPara.Range.InsertXML($checkXml)
}

# ---------------------------------------------------------------------------
# Change 2: paragraph "Record expected initial value of every sensor"
# (numId=6) loses the eastAsia rPr that was on its paragraph mark, and is
# followed by four new paragraphs:
#   - an empty paragraph
#   - a date line "202108113" (split "2" + "02108113", first run hinted eastAsia)
#   - a new numId=6 bullet paragraph with the "hip joint velocity feedback" text
#   - a trailing empty paragraph whose mark carries an eastAsia rPr
# ---------------------------------------------------------------------------
$recordPara = $d.Paragraphs.Item(6)
if ($recordPara.Range.Text.StartsWith("Record expected")) {
    $recordXml = "<w:p $W w:rsidR='00D73537' w:rsidRDefault='00574CE4' w:rsidP='00574CE4'>" +
        "<w:pPr>" +
          "<w:pStyle w:val='a7'/>" +
          "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr>" +
          "<w:ind w:firstLineChars='0'/>" +
        "</w:pPr>" +
        "<w:r><w:t>Record expected initial value of every sensor</w:t></w:r>" +
      "</w:p>"
    $recordXml += "<w:p $W/>"
    $recordXml += "<w:p $W>" +
        "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>2</w:t></w:r>" +
        "<w:r><w:t>02108113</w:t></w:r>" +
      "</w:p>"
    $recordXml += "<w:p $W>" +
        "<w:pPr>" +
          "<w:pStyle w:val='a7'/>" +
          "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='6'/></w:numPr>" +
          "<w:ind w:firstLineChars='0'/>" +
        "</w:pPr>" +
        "<w:r><w:t>The hip joint velocity feedback is very essential</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>, so the </w:t></w:r>" +
        "<w:r><w:t>velocity feedback may be obtained from motor driver or potentiometer feedback differentiation. Should be selected</w:t></w:r>" +
      "</w:p>"
    $recordXml += "<w:p $W>" +
        "<w:pPr><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr></w:pPr>" +
      "</w:p>"
    $recordPara.Range.InsertXML($recordXml)
}
